# Update "gh-pages" output generated at 456a3b4
# Applies numeric "want-to-go count" bumps across the four sheets and
# reshuffles the 展览 (Exhibition) sheet to insert two newly scraped
# events.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: 展览 (Exhibition)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Small "want to go" (column F) bumps for existing rows.
$ws1.Range("F2").Value2  = 223
$ws1.Range("F3").Value2  = 2507
$ws1.Range("F5").Value2  = 1960
$ws1.Range("F6").Value2  = 1325
$ws1.Range("F12").Value2 = 1818
$ws1.Range("F14").Value2 = 1875
$ws1.Range("F18").Value2 = 528
$ws1.Range("F21").Value2 = 25
$ws1.Range("F22").Value2 = 18
$ws1.Range("F23").Value2 = 2416
$ws1.Range("F24").Value2 = 460
$ws1.Range("F26").Value2 = 1047
$ws1.Range("F27").Value2 = 4617
$ws1.Range("F28").Value2 = 111

# Two brand-new events were scraped: one that now sorts in as row 31
# ("SunShine") and one that sorts in as row 35 ("旅行盛宴 4.0"). Both
# push the rows that already existed further down the sheet, so insert
# two blank rows at the right spots first.
$ws1.Rows.Item(31).Insert()
$ws1.Rows.Item(35).Insert()

# Give the new A-column cells the same style as the rest of the index
# column (bold/centered/bordered, style index 1) by copying the format
# from a neighbouring cell.
$ws1.Range("A30").Copy()
$ws1.Range("A31").PasteSpecial(-4122) | Out-Null
$ws1.Range("A30").Copy()
$ws1.Range("A35").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Helper: write a "YYYY-MM-DD"-shaped string into column B without
# Excel's autodetection turning it into a date serial number - force
# the cell to Text first, write the value, then drop back to the
# workbook's default "Normal" style so no stray numFmt lingers.
function Set-TextDate($range, $text) {
    $range.NumberFormat = "@"
    $range.Value2 = $text
    $range.Style = "Normal"
}

$ws1.Range("A31").Value2 = 30
Set-TextDate $ws1.Range("B31") "2024-11-30"
$ws1.Range("C31").Value2 = "上海·SunShine跨次元动漫游戏嘉年华2.0"
$ws1.Range("D31").Value2 = "景泰路465号 上海园艺体验中心"
$ws1.Range("E31").Value2 = "2024.11.30 10:30-12.01 17:00"
$ws1.Range("F31").Value2 = 0
$ws1.Range("G31").Value2 = 69
$ws1.Range("H31").Value2 = "https://show.bilibili.com/platform/detail.html?id=93445"
$ws1.Range("I31").Value2 = "//i2.hdslb.com/bfs/openplatform/202410/PkBLAxyI1728882644725.jpeg"

$ws1.Range("A35").Value2 = 34
Set-TextDate $ws1.Range("B35") "2024-12-21"
$ws1.Range("C35").Value2 = "上海·旅行盛宴次元综合同人动漫节4.0·一周年庆"
$ws1.Range("D35").Value2 = "景泰路465号 上海园艺体验中心"
$ws1.Range("E35").Value2 = "2024.12.21 10:00-12.22 17:00"
$ws1.Range("F35").Value2 = 2
$ws1.Range("G35").Value2 = 69
$ws1.Range("H35").Value2 = "https://show.bilibili.com/platform/detail.html?id=93447"
$ws1.Range("I35").Value2 = "//i1.hdslb.com/bfs/openplatform/202410/9uMCI2Ac1728894507590.jpeg"

# Column A is a plain 0-based row index (row N -> value N-1). The row
# Inserts above correctly shifted columns B:I but left column A's old
# values behind (A is the same every time something is inserted above
# it, values don't actually travel with "their" event) - rewrite A for
# every row from the first new row through the end of the sheet so the
# sequence is contiguous again.
for ($r = 32; $r -le 37; $r++) {
    $ws1.Cells.Item($r, 1).Value2 = $r - 1
}

# ---------------------------------------------------------------------
# Sheet 2: 演出 (Shows)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F4").Value2  = 12
$ws2.Range("F10").Value2 = 4
$ws2.Range("F19").Value2 = 190
$ws2.Range("F32").Value2 = 485
$ws2.Range("F42").Value2 = 104

# ---------------------------------------------------------------------
# Sheet 3: 本地生活 (Local life)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F9").Value2  = 3128
$ws3.Range("F10").Value2 = 631
$ws3.Range("F11").Value2 = 901
$ws3.Range("F14").Value2 = 65
$ws3.Range("F15").Value2 = 18
$ws3.Range("F16").Value2 = 322

# ---------------------------------------------------------------------
# Sheet 4: 全部类型 (All types)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F4").Value2  = 2507
$ws4.Range("F6").Value2  = 3128
$ws4.Range("F7").Value2  = 631
$ws4.Range("F8").Value2  = 901
$ws4.Range("F9").Value2  = 1960
$ws4.Range("F12").Value2 = 65
$ws4.Range("F13").Value2 = 65
$ws4.Range("F14").Value2 = 1325
$ws4.Range("F17").Value2 = 18
$ws4.Range("F18").Value2 = 1818
$ws4.Range("F21").Value2 = 1875
$ws4.Range("F23").Value2 = 529
$ws4.Range("F26").Value2 = 190
$ws4.Range("F28").Value2 = 25
$ws4.Range("F29").Value2 = 18
$ws4.Range("F31").Value2 = 2416
$ws4.Range("F32").Value2 = 460
$ws4.Range("F35").Value2 = 1047
$ws4.Range("F37").Value2 = 322
$ws4.Range("F40").Value2 = 111
$ws4.Range("F41").Value2 = 485
$ws4.Range("F46").Value2 = 176
$ws4.Range("F47").Value2 = 104

Write-Host "edit complete"
